# Automatische test-sync: 2025-07-29 21:39:50
#
# This script appends a new test-mail log entry (row 8) to the "Logs"
# worksheet, adds the corresponding summary row (row 5) to the
# "Dashboard" worksheet, extends the conditional formatting ranges on
# "Logs" so they keep covering the full data range, and updates the
# embedded bar chart's category/value series so it plots the new
# Dashboard row as well.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append row 8 with the new test mail entry.
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A8").Value = "Hebben we EcoPro-700 nog op voorraad?"
$logs.Range("B8").Value = "mailmind.test@zohomail.eu"
$logs.Range("C8").Value = "Testmail #6: Hebben we EcoPro-700 nog op voorraad?"
$logs.Range("D8").Value = "Productinformatie"
$logs.Range("E8").Value = "Beste afzender,`nHartelijk dank voor uw interesse in onze EcoPro-700. Op dit moment hebben we nog voldoende EcoPro-700 op voorraad. U kunt deze direct bestellen via onze website of neem contact met ons op als u meer informatie wenst.`nMet vriendelijke groet,`n[Naam]  `nE-mailassistent  `n[Bedrijfsnaam]"
$logs.Range("F8").Value = "2025-07-29 21:39:42"
$logs.Range("G8").Value = "Ja"
$logs.Range("H8").Value = "Nee"
$logs.Range("I8").Value = "Ja"
$logs.Range("J8").Value = "Nee"

# Setting the multi-line answer above makes the engine auto-expand the
# row height; restore the (default) standard row height/AutoFit so the
# row stays un-customized, matching the other data rows.
$logs.Rows.Item(8).AutoFit()

# ---------------------------------------------------------------------
# 2. Extend the conditional formatting ranges from row 7 to row 8 so
#    that they keep applying to the full data range.
# ---------------------------------------------------------------------
$logs.Range("D2:D7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D8"))
$logs.Range("G2:G7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G8"))
$logs.Range("H2:H7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H8"))
$logs.Range("I2:I7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I8"))
$logs.Range("J2:J7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J8"))

# ---------------------------------------------------------------------
# 3. Dashboard sheet: append row 5 with the "Productinformatie" count.
# ---------------------------------------------------------------------
$dashboard = $wb.Worksheets.Item("Dashboard")

$dashboard.Range("A5").Value = "Productinformatie"
$dashboard.Range("B5").Value = 1

# ---------------------------------------------------------------------
# 4. Update the bar chart series so the category/value references
#    cover the new Dashboard row (A2:A5 / B2:B5 instead of A2:A4 /
#    B2:B4). The series name reference stays untouched.
# ---------------------------------------------------------------------
$chartObj = $dashboard.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES('Dashboard'!B1,'Dashboard'!`$A`$2:`$A`$5,'Dashboard'!`$B`$2:`$B`$5,1)"
